$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target dataset (Regiao, Variavel, Ano, Valor) for rows 2..57.
# $null in the 4th slot means the Valor cell is blank (no <v> in the source).
$data = @(
    @('Brasil', 'Roubo a instituição financeira', '31/12/2016', 1.503707104007567),
    @('Brasil', 'Roubo de carga', '31/12/2016', 11.33726594058217),
    @('Brasil', 'Roubo de veículo', '31/12/2016', 295.7215713339049),
    @('Brasil', 'Furto de veículo', '31/12/2016', 302.0852233652141),
    @('Brasil', 'Roubo a instituição financeira', '31/12/2017', 2.150349989335505),
    @('Brasil', 'Roubo de carga', '31/12/2017', 13.42443379025678),
    @('Brasil', 'Roubo de veículo', '31/12/2017', 292.5194070621019),
    @('Brasil', 'Furto de veículo', '31/12/2017', 276.9874770058215),
    @('Brasil', 'Roubo a instituição financeira', '31/12/2018', 1.213026853601723),
    @('Brasil', 'Roubo de carga', '31/12/2018', 10.96543073447625),
    @('Brasil', 'Roubo de veículo', '31/12/2018', 245.1617376924052),
    @('Brasil', 'Furto de veículo', '31/12/2018', 243.4336388660364),
    @('Brasil', 'Roubo a instituição financeira', '31/12/2019', 0.7110696644275744),
    @('Brasil', 'Roubo de carga', '31/12/2019', 8.436185683820611),
    @('Brasil', 'Roubo de veículo', '31/12/2019', 180.4868330798366),
    @('Brasil', 'Furto de veículo', '31/12/2019', 218.2787271480123),
    @('Brasil', 'Roubo a instituição financeira', '31/12/2020', 0.5040368218516562),
    @('Brasil', 'Roubo de carga', '31/12/2020', 6.293685499770147),
    @('Brasil', 'Roubo de veículo', '31/12/2020', 133.7588438583084),
    @('Brasil', 'Furto de veículo', '31/12/2020', 165.687945739034),
    @('Brasil', 'Roubo a instituição financeira', '31/12/2021', 0.4870619219114404),
    @('Brasil', 'Roubo de carga', '31/12/2021', 6.768525769265649),
    @('Brasil', 'Roubo de veículo', '31/12/2021', 128.083453577476),
    @('Brasil', 'Furto de veículo', '31/12/2021', 172.2524822814674),
    @('Brasil', 'Roubo a instituição financeira', '31/12/2022', 0.3802159882420602),
    @('Brasil', 'Roubo de carga', '31/12/2022', 6.467526599564192),
    @('Brasil', 'Roubo de veículo', '31/12/2022', 128.5514751260922),
    @('Brasil', 'Furto de veículo', '31/12/2022', 195.6634690836586),
    @('Sergipe', 'Roubo a instituição financeira', '31/12/2016', $null),
    @('Sergipe', 'Roubo de carga', '31/12/2016', $null),
    @('Sergipe', 'Roubo de veículo', '31/12/2016', 378.1975589066652),
    @('Sergipe', 'Furto de veículo', '31/12/2016', 142.599079587759),
    @('Sergipe', 'Roubo a instituição financeira', '31/12/2017', $null),
    @('Sergipe', 'Roubo de carga', '31/12/2017', $null),
    @('Sergipe', 'Roubo de veículo', '31/12/2017', 380.4861857878965),
    @('Sergipe', 'Furto de veículo', '31/12/2017', 116.9792635028746),
    @('Sergipe', 'Roubo a instituição financeira', '31/12/2018', 0.4622496147919877),
    @('Sergipe', 'Roubo de carga', '31/12/2018', $null),
    @('Sergipe', 'Roubo de veículo', '31/12/2018', 341.6711981149175),
    @('Sergipe', 'Furto de veículo', '31/12/2018', 100.4686812190891),
    @('Sergipe', 'Roubo a instituição financeira', '31/12/2019', 0.1592356687898089),
    @('Sergipe', 'Roubo de carga', '31/12/2019', 0.04350292513668619),
    @('Sergipe', 'Roubo de veículo', '31/12/2019', 239.9654747347979),
    @('Sergipe', 'Furto de veículo', '31/12/2019', 109.2555985743447),
    @('Sergipe', 'Roubo a instituição financeira', '31/12/2020', 0.1636661211129296),
    @('Sergipe', 'Roubo de carga', '31/12/2020', $null),
    @('Sergipe', 'Roubo de veículo', '31/12/2020', 301.6173092963617),
    @('Sergipe', 'Furto de veículo', '31/12/2020', 94.16578436035383),
    @('Sergipe', 'Roubo a instituição financeira', '31/12/2021', 0.1709401709401709),
    @('Sergipe', 'Roubo de carga', '31/12/2021', 0.09100162646448628),
    @('Sergipe', 'Roubo de veículo', '31/12/2021', 273.3500348713431),
    @('Sergipe', 'Furto de veículo', '31/12/2021', 253.2786583752849),
    @('Sergipe', 'Roubo a instituição financeira', '31/12/2022', $null),
    @('Sergipe', 'Roubo de carga', '31/12/2022', 0.09051584072470603),
    @('Sergipe', 'Roubo de veículo', '31/12/2022', 209.8330592866558),
    @('Sergipe', 'Furto de veículo', '31/12/2022', 287.7139022898047)
)

$totalRows = $data.Length
$lastRowBefore = $ws.UsedRange.Rows.Count

# Clear any existing data rows beyond the header first (in case the new
# table is shorter or longer than the old one), then write the new table.
if ($lastRowBefore -gt 1) {
    $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRowBefore, 4)).ClearContents()
}

for ($i = 0; $i -lt $totalRows; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($null -eq $row[3]) {
        $ws.Cells.Item($r, 4).ClearContents()
    } else {
        $ws.Cells.Item($r, 4).Value = $row[3]
    }
}
